$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Lot No N" labels to "Lot No. N" across the three test-kit blocks
$ws.Range("T1").Value = "Lot No. 1"
$ws.Range("Z1").Value = "Lot No. 2"
$ws.Range("AF1").Value = "Lot No. 3"

# Move the active selection to AF1 (matches the saved view state in the
# workbook after the edit)
$null = $ws.Range("AF1").Select()
